$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so numeric-looking
# strings (e.g. "0.9579", "20.834.88") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '20.834.88'
$ws.Range('E2').Value = '  +2.55%  '
$ws.Range('D3').Value = '1.520.21'
$ws.Range('E3').Value = '  +4.45%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '0.9579'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').Value = '279.72'
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('D7').Value = '0.3588'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('D8').Value = '0.3131'
$ws.Range('E8').Value = '  +2.14%  '
$ws.Range('E9').Value = '  +7.30%  '
$ws.Range('D10').Value = '40.09'
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').Value = '0.06741'
$ws.Range('E11').Value = '  +3.15%  '
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '18.63'
$ws.Range('E13').Value = '  +5.19%  '
$ws.Range('D14').Value = '5.603'
$ws.Range('E14').Value = '  +4.37%  '
$ws.Range('D15').Value = '6.272'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('D16').Value = '0.9586'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '0.00001031'
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').Value = '1.516.23'
$ws.Range('E18').Value = '  +4.13%  '
$ws.Range('D19').Value = '0.06041'
$ws.Range('E19').Value = '  +5.54%  '
$ws.Range('D20').Value = '70.57'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Value = '5.617'
$ws.Range('E21').Value = '  +3.92%  '
$ws.Range('E22').Value = '  +3.62%  '
$ws.Range('D24').Value = '2.293'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('D25').Value = '20.869.43'
$ws.Range('D26').Value = '146.44'
$ws.Range('E26').Value = '  +4.33%  '
$ws.Range('D27').Value = '2.164'
$ws.Range('E27').Value = '  +3.77%  '
$ws.Range('D28').Value = '17.41'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('D29').Value = '1.677.03'
$ws.Range('E29').Value = '  +4.30%  '
$ws.Range('D30').Value = '116.50'
$ws.Range('E30').Value = '  +4.71%  '
$ws.Range('D31').Value = '4.006'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.127'
$ws.Range('E32').Value = '  +5.90%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.8376'
$ws.Range('E33').Value = '  +6.19%  '
$ws.Range('D34').Value = '0.08003'
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').Value = '1.210'
$ws.Range('E35').Value = '  +7.46%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.468'
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('D37').Value = '4.857'
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('D38').Value = '0.05797'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('D39').Value = '0.02077'
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('D41').Value = '0.9590'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').Value = '0.1885'
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('D43').Value = '7.543'
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('D44').Value = '0.5334'
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('D45').Value = '3.544'
$ws.Range('E45').Value = '  +1.92%  '
$ws.Range('D46').Value = '12.30'
$ws.Range('E46').Value = '  +3.38%  '
$ws.Range('D47').Value = '120.66'
$ws.Range('E47').Value = '  +3.27%  '
$ws.Range('D48').Value = '0.5330'
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('D49').Value = '1.852'
$ws.Range('E49').Value = '  +6.45%  '
$ws.Range('D50').Value = '0.06515'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').Value = '0.9876'
$ws.Range('E51').Value = '  +0.40%  '

# Restore default (un-styled) formatting to match the original workbook look.
$ws.Range("D2:E51").Style = "Normal"
